$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Apply the border/shading format used by the LAST data row
#    (previously row 26) onto what will become the new last data
#    row (row 19), before we touch any values.
# ------------------------------------------------------------------
$ws.Range("B26:J26").Copy() | Out-Null
$ws.Range("B19:J19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ------------------------------------------------------------------
# 2) Overwrite the 4 data rows (16-19) with the new "Base de datos"
#    worker records.
# ------------------------------------------------------------------

# Row 16: CC 1047475389 - CARLOS FERNANDO OZUNA CORTINA - period 2008
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047475389"
$ws.Range("D16").Value = "CARLOS FERNANDO OZUNA CORTINA"
$ws.Range("E16").Value = "2008"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 2064203

# Row 17: CC 1235044752 - JOSE CARLOS RIOS MARQUEZ - period 2008
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1235044752"
$ws.Range("D17").Value = "JOSE CARLOS RIOS MARQUEZ"
$ws.Range("E17").Value = "2008"
$ws.Range("F17").Value = 26666
$ws.Range("G17").Value = 5511848

# Row 18: CE 20394544 - LUZ ELENA UTRIA ORTIZ - period 2210
$ws.Range("B18").Value = "CE"
$ws.Range("C18").Value = "20394544"
$ws.Range("D18").Value = "LUZ ELENA UTRIA ORTIZ"
$ws.Range("E18").Value = "2210"
$ws.Range("F18").Value = 9333
$ws.Range("G18").Value = 1000000

# Row 19 (last data row, bottom-border style): CC 1143352669 - CARLOS MARIO FUENTES MAYA - period 2211
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143352669"
$ws.Range("D19").Value = "CARLOS MARIO FUENTES MAYA"
$ws.Range("E19").Value = "2211"
$ws.Range("F19").Value = 105620
$ws.Range("G19").Value = 2640500

# ------------------------------------------------------------------
# 3) Remove the now-obsolete trailing data rows (20-26). This shifts
#    the signature block (previously rows 31-32) up to rows 24-25.
# ------------------------------------------------------------------
$ws.Rows("20:26").Delete()

# ------------------------------------------------------------------
# 4) Update the summary figures above the table.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 181619   # VALOR MORA total
$ws.Range("C13").Value = 4        # Cant. Trabajadores
$ws.Range("F13").Value = 3        # Cant. Periodos
